$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address, new value, and whether the value needs to be
# forced to Text so Excel does not auto-convert numeric-looking strings
# (e.g. "1.003") into real numbers. Non-numeric-looking strings (coin names,
# links, "xx.xx.xx" prices, and the "  +n.nn%  " volume strings) are assigned
# directly since Excel already stores them as text.
$updates = @(
    [PSCustomObject]@{ Addr = 'D2'; Value = '29.248.56'; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E2'; Value = '  -0.38%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D3'; Value = '1.829.36'; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E3'; Value = '  -0.54%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D4'; Value = '1.003'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E4'; Value = '  +0.34%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D5'; Value = '235.64'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E5'; Value = '  -1.31%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D6'; Value = '0.6025'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E6'; Value = '  -3.75%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D7'; Value = '1.005'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'D8'; Value = '0.07037'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E8'; Value = '  -4.93%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D9'; Value = '0.2792'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E9'; Value = '  -3.28%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D10'; Value = '23.55'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E10'; Value = '  -5.12%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D12'; Value = '1.822.57'; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E12'; Value = '  -0.81%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D13'; Value = '4.789'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E13'; Value = '  -3.51%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D14'; Value = '0.6275'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E14'; Value = '  -6.75%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D15'; Value = '0.000009766'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E15'; Value = '  -4.36%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D16'; Value = '78.97'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E16'; Value = '  -3.25%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D17'; Value = '29.263.25'; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E17'; Value = '  -0.62%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D18'; Value = '5.817'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E18'; Value = '  -6.18%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D19'; Value = '224.12'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E19'; Value = '  -3.68%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E20'; Value = '  +0.32%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D21'; Value = '11.68'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E21'; Value = '  -5.04%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D22'; Value = '6.997'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E22'; Value = '  -3.91%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D23'; Value = '1.004'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E23'; Value = '  +0.23%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D24'; Value = '156.24'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E24'; Value = '  -1.03%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'B25'; Value = 'Stellar'; ForceText = $false }
    [PSCustomObject]@{ Addr = 'C25'; Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D25'; Value = '0.1300'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E25'; Value = '  -3.11%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'B26'; Value = 'Cosmos'; ForceText = $false }
    [PSCustomObject]@{ Addr = 'C26'; Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D26'; Value = '7.982'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E26'; Value = '  -5.85%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D27'; Value = '16.60'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E27'; Value = '  -3.91%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D28'; Value = '0.06642'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E28'; Value = '  -8.57%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D29'; Value = '1.469'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E29'; Value = '  -0.45%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D30'; Value = '1.447'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E30'; Value = '  -1.89%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D31'; Value = '3.842'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E31'; Value = '  -4.43%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E32'; Value = '  -5.92%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E33'; Value = '  -2.94%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D34'; Value = '1.723'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E34'; Value = '  -4.93%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D35'; Value = '0.6454'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E35'; Value = '  -7.29%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D36'; Value = '2.546'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E36'; Value = '  -1.00%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E37'; Value = '  -2.69%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D38'; Value = '1.211.64'; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E38'; Value = '  -1.74%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D39'; Value = '0.01757'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E39'; Value = '  -4.37%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D40'; Value = '6.523'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E40'; Value = '  -5.40%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D41'; Value = '0.9017'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E41'; Value = '  -4.46%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E42'; Value = '  +0.41%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D43'; Value = '1.988.78'; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E43'; Value = '  -1.43%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D44'; Value = '100.53'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E44'; Value = '  -0.13%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D45'; Value = '62.55'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E45'; Value = '  -4.16%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'E46'; Value = '  -0.80%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D47'; Value = '8.540'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E47'; Value = '  -3.48%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D48'; Value = '1.582'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E48'; Value = '  -7.02%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D49'; Value = '0.4552'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E49'; Value = '  -0.47%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D50'; Value = '0.05501'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E50'; Value = '  -2.77%  '; ForceText = $false }
    [PSCustomObject]@{ Addr = 'D51'; Value = '6.388'; ForceText = $true }
    [PSCustomObject]@{ Addr = 'E51'; Value = '  -7.80%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $cell = $ws.Range($u.Addr)
    if ($u.ForceText) {
        $cell.NumberFormat = "@"
        $cell.Value = $u.Value
        $cell.Style = "Normal"
    } else {
        $cell.Value = $u.Value
    }
}
